$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column U: header "float" / "radius", plus per-row radius values
$ws.Range("U4").Value = "float"
$ws.Range("U5").Value = "radius"

$ws.Range("U6").Value = "'12"
$ws.Range("U7").Value = "'14"
$ws.Range("U8").Value = "'18"
$ws.Range("U9").Value = "'20"
$ws.Range("U10").Value = "'24"
$ws.Range("U11").Value = "'10"

# Update projectile/impact sprite paths to new "ui/assets/" prefix
$ws.Range("M6").Value = "ui/assets/fx/projectiles/spittle.png"
$ws.Range("N6").Value = "ui/assets/fx/impact/slime.png"

$ws.Range("M7").Value = "ui/assets/fx/projectiles/choir_note.png"
$ws.Range("N7").Value = "ui/assets/fx/impact/chorus.png"

$ws.Range("M8").Value = "ui/assets/fx/projectiles/howl_wave.png"
$ws.Range("N8").Value = "ui/assets/fx/impact/howl.png"

$ws.Range("M9").Value = "ui/assets/fx/projectiles/null_beam.png"
$ws.Range("N9").Value = "ui/assets/fx/impact/null_burn.png"

$ws.Range("M10").Value = "ui/assets/fx/projectiles/dredger_slam.png"
$ws.Range("N10").Value = "ui/assets/fx/impact/dredger_slam.png"

$ws.Range("M11").Value = "ui/assets/fx/projectiles/fragment_dart.png"
$ws.Range("N11").Value = "ui/assets/fx/impact/fragment_spark.png"
